# Trade #64 closed at 2026-02-17 08:48:48 - unknown UNKNOWN +0.000%
#
# Applies:
#  - Summary sheet metric updates (new trade pushes totals)
#  - Strategy Status sheet (MarketMaking row) updates
#  - All Trades sheet: append the new closed trade as row 65
#  - MarketMaking sheet: append the same new closed trade as row 65

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.59   # Current Capital
$summary.Range("B4").Value = -0.41     # Total P&L $
$summary.Range("B5").Value = -0.13     # Total P&L %
$summary.Range("B6").Value = 64        # Total Trades
$summary.Range("B8").Value = 27        # Losing Trades
$summary.Range("B9").Value = 39.06     # Win Rate %

# ---------------------------------------------------------------------
# Sheet 2: Strategy Status (MarketMaking is row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.59      # Capital
$status.Range("D4").Value = 64         # Trades
$status.Range("E4").Value = -0.41      # P&L $
$status.Range("F4").Value = -0.41      # P&L %
$status.Range("G4").Value = 39.06      # Win Rate %

# ---------------------------------------------------------------------
# Helper: append trade #64 as a new last row (row 65) on a trades sheet.
# Copy the previous last row first so text-valued cells (dates/times/
# labels) keep their literal-text representation instead of Excel
# auto-converting the new values into dates/numbers, then overwrite the
# cells that actually changed.
# ---------------------------------------------------------------------
function Add-Trade64Row($ws) {
    $ws.Range("A64:Q64").Copy($ws.Range("A65:Q65"))

    $ws.Cells.Item(65, 1).Value = 64                 # Trade #
    # B65 Date stays "2026-02-17" (copied)
    $ws.Cells.Item(65, 3).Value = "08:48:41"          # Time
    # D65 Strategy stays "MarketMaking" (copied)
    $ws.Cells.Item(65, 5).Value = "DOWN"              # Side
    $ws.Cells.Item(65, 6).Value = 0.18                # Entry Price
    $ws.Cells.Item(65, 7).Value = 0.12                # Exit Price
    # H65 Status stays "CLOSED" (copied)
    $ws.Cells.Item(65, 9).Value = -33.3333            # P&L %
    $ws.Cells.Item(65, 10).Value = -0.06              # P&L $
    $ws.Cells.Item(65, 11).Value = 99.59              # Capital After
    # L65 Entry Slippage stays 0 (copied)
    # M65 Exit Slippage stays 0 (copied)
    # N65 Confidence stays 0.6 (copied)
    # O65 Entry Reason stays "Normal spread capture: 19600 bps" (copied)
    # P65 Exit Reason stays "early_exit" (copied)
    # Q65 Duration (min) stays 0.14 (copied)
}

# ---------------------------------------------------------------------
# Sheet 3: All Trades
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade64Row $allTrades

# ---------------------------------------------------------------------
# Sheet 4: MarketMaking
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade64Row $marketMaking
